$d = $word.ActiveDocument

# 1. "CAG/$" + "NUMERO" + "$/21" -> "CAG/$NUMERO$/21" (merge runs, drop proofErr marks)
$d.Content.Find.Execute("CAG/`$NUMERO`$/21", $true, $false, $false, $false, $false, $true, 1, $false, "CAG/`$NUMERO`$/21", 2)

# 2. "de acuerdo al" -> "de acuerdo con el" (drop proofErr marks around it)
$d.Content.Find.Execute("de acuerdo al", $true, $false, $false, $false, $false, $true, 1, $false, "de acuerdo con el", 2)

# 3. "Año del Legado de Fray Antonio " + "Alcalde" + " en Guadalajara" -> merge single run
$d.Content.Find.Execute("Año del Legado de Fray Antonio Alcalde en Guadalajara", $true, $false, $false, $false, $false, $true, 1, $false, "Año del Legado de Fray Antonio Alcalde en Guadalajara", 2)

# 4. Footer text merge
$d.Content.Find.Execute("Calle Sierra Nevada No. 910, Col. Independencia, C.P. 44340, Guadalajara, Jal., México,  Tel (52) 33 1058 5249. y fax:  (52) 33 1058 5200   ext. 337", $true, $false, $false, $false, $false, $true, 1, $false, "Calle Sierra Nevada No. 910, Col. Independencia, C.P. 44340, Guadalajara, Jal., México,  Tel (52) 33 1058 5249. y fax:  (52) 33 1058 5200   ext. 337", 2)

Write-Output "done"
